$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 16, shifting existing rows 16-35 down to 17-36.
$ws.Rows.Item(16).Insert()

# Populate the newly inserted row 16 with its data.
$ws.Range("A16").Value = 5
$ws.Range("B16").Value = "Macroferia Regional de Talca"
$ws.Range("C16").Value = "Maule"
$ws.Range("D16").Value = 44645
$ws.Range("E16").Value = 7
$ws.Range("F16").Value = "Fruta"
$ws.Range("G16").Value = 100107
$ws.Range("H16").Value = "Otros"
$ws.Range("I16").Value = 100107011
$ws.Range("J16").Value = "Tuna"
$ws.Range("K16").Value = "Sin especificar"
$ws.Range("L16").Value = "Primera"
$ws.Range("M16").Value = 200
$ws.Range("N16").Value = 16000
$ws.Range("O16").Value = 16000
$ws.Range("P16").Value = 16000
$ws.Range("Q16").Value = "$/caja 18 kilos"
$ws.Range("R16").Value = "Provincia de Limarí"
$ws.Range("S16").Value = 889
$ws.Range("T16").Value = 18
